$wb = $excel.ActiveWorkbook

# The workbook currently ends with the "Greece" sheet, which acts as the
# template for the three new country sheets being added (Netherlands,
# Austria, Denmark) for the Austria-market test data work.
$template = $wb.Worksheets.Item("Greece")

# --- Create "Netherlands" (copied right after Greece) ---
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy([System.Reflection.Missing]::Value, $afterSheet)
$netherlands = $wb.Worksheets.Item($wb.Worksheets.Count)
$netherlands.Name = "Netherlands"

# --- Create "Austria" (copied right after Netherlands) ---
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy([System.Reflection.Missing]::Value, $afterSheet)
$austria = $wb.Worksheets.Item($wb.Worksheets.Count)
$austria.Name = "Austria"

# --- Create "Denmark" (copied right after Austria) ---
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy([System.Reflection.Missing]::Value, $afterSheet)
$denmark = $wb.Worksheets.Item($wb.Worksheets.Count)
$denmark.Name = "Denmark"

# --- Fill in the market-specific data ---
# Netherlands
$netherlands.Range("B4").Value = "NGC-4330/T2199"
$netherlands.Range("B2").Value = "Netherlands Market"
$netherlands.Rows.Item(2).RowHeight = 28.8

# Austria market name first ...
$austria.Range("B2").Value = "Austria Market"

# Denmark
$denmark.Range("B4").Value = "NGC-2913/T2798"
$denmark.Range("B2").Value = "Denmark Market"

# ... then Austria's user-story reference ...
$austria.Range("B4").Value = "NGC-3817/T2306"

# Austria gets an extra "Fire Brigade Panel" row inserted just above the
# trailing "Wg" / "Miscellaneous" rows, pushing them down by one row.
$austria.Rows.Item(9).Insert()
$austria.Range("A8").Copy()
$austria.Range("A9").PasteSpecial(-4122)
$austria.Range("A9").Value = "Fire Brigade Panel"

# Austria ends up the active/selected sheet.
$austria.Activate()

Write-Output "Added Netherlands, Austria and Denmark sheets"
